$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 23.9300000000003
$ws.Range("H2").Value = [double]"1.154978439141905e-16"
$ws.Range("K2").Value = 60.5171714126482
$ws.Range("L2").Value = "[56.36013662065184, 64.67420620464455]"
$ws.Range("O2").Value = 1.628973968528041
$ws.Range("P2").Value = "[1.566079220708426, 1.6918687163476562]"
$ws.Range("S2").Value = 56.17493854105904
$ws.Range("T2").Value = "[53.43450829573093, 58.915368786387155]"
$ws.Range("W2").Value = 17.72592592592615
$ws.Range("X2").Value = 17.48638638638661
$ws.Range("Y2").Value = 17.96546546546569

# Row 3
$ws.Range("E3").Value = 23.88000000000029
$ws.Range("H3").Value = [double]"1.154978439141905e-16"
$ws.Range("K3").Value = 55.51302827476059
$ws.Range("L3").Value = "[50.01468705530349, 61.01136949421768]"
$ws.Range("O3").Value = 1.50318447288881
$ws.Range("P3").Value = "[1.4025528763774249, 1.603816069400195]"
$ws.Range("S3").Value = 52.10518401442899
$ws.Range("T3").Value = "[48.53789929944419, 55.67246872941378]"
$ws.Range("W3").Value = 18.16696696696719
$ws.Range("X3").Value = 17.78450450450472
$ws.Range("Y3").Value = 18.54942942942966
